$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the raw measurement data in columns C (Sklon mV) and G (Sklon mBar):
# the values had been entered one decimal place too large (e.g. 16.081
# instead of 1.6081) and C12 had a typo that should match the corrected
# C11 value. AVERAGE formulas in row 5 recalculate automatically.

$ws.Range("C7").Value = 1.6081000000000001
$ws.Range("C8").Value = 1.6288
$ws.Range("C9").Value = 1.6408
$ws.Range("C10").Value = 1.6297999999999999
$ws.Range("C11").Value = 1.6015999999999999
$ws.Range("C12").Value = 1.6015999999999999

$ws.Range("G7").Value = 0.26336999999999999
$ws.Range("G8").Value = 0.26433000000000001
$ws.Range("G9").Value = 0.26977000000000001
$ws.Range("G10").Value = 0.26734999999999998
$ws.Range("G11").Value = 0.26001999999999997
$ws.Range("G12").Value = 0.26157999999999998

# Update the active selection to match where the author ended up (C13).
$ws.Range("C13").Select()
